# Apply the "Generate Report for Archive" update:
#  1. Change the localization status text "Ready for handoff" -> "In Translation"
#     on every sheet/cell that currently shows it.
#  2. Narrow the corresponding "Status" columns (Overview!E:F, zh-cn!C, de-de!C)
#     to reflect the shorter text, matching the new auto-fit column width.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

# Columns E and F narrow from ~17.22 to ~13.41 raw width units (status text shrank)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet ------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
